$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "disponible"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4108
$ws.Range("F1").NumberFormat = "@"

# Fill column F (rows 2-13 and 15-19, skipping 14) with "1"
$rows = 2,3,4,5,6,7,8,9,10,11,12,13,15,16,17,18,19
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "1"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.NumberFormat = "@"
}

# Leave the selection where the user last clicked
$ws.Range("I10").Select() | Out-Null
